# Append new trading log rows (44-47) to the active worksheet,
# mirroring the data appended by the latest trading run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44: TRADING_ATTEMPT for SOL
$ws.Cells.Item(44, 1).Value = "2025-09-25T01:26:48.157640"
$ws.Cells.Item(44, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(44, 3).Value = "SOL"
$ws.Cells.Item(44, 4).Value = "UNKNOWN"
$ws.Cells.Item(44, 5).Value = 210.841151833073
$ws.Cells.Item(44, 11).Value = "ATTEMPT"
$ws.Cells.Item(44, 12).Value = "Attempting trade 1/2"

# Row 45: POSITION_OPENED for SOL
$ws.Cells.Item(45, 1).Value = "2025-09-25T01:26:49.937087"
$ws.Cells.Item(45, 2).Value = "POSITION_OPENED"
$ws.Cells.Item(45, 3).Value = "SOL"
$ws.Cells.Item(45, 4).Value = "UNKNOWN"
$ws.Cells.Item(45, 5).Value = 210.841151833073
$ws.Cells.Item(45, 6).Value = 2400
$ws.Cells.Item(45, 7).Value = 20
$ws.Cells.Item(45, 8).Value = 0.4434054849163298
$ws.Cells.Item(45, 11).Value = "SUCCESS"

# Row 46: TRADING_ATTEMPT for SUI
$ws.Cells.Item(46, 1).Value = "2025-09-25T01:26:49.962191"
$ws.Cells.Item(46, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(46, 3).Value = "SUI"
$ws.Cells.Item(46, 4).Value = "UNKNOWN"
$ws.Cells.Item(46, 5).Value = 3.348909889983477
$ws.Cells.Item(46, 11).Value = "ATTEMPT"
$ws.Cells.Item(46, 12).Value = "Attempting trade 2/2"

# Row 47: POSITION_OPENED for SUI
$ws.Cells.Item(47, 1).Value = "2025-09-25T01:26:51.589220"
$ws.Cells.Item(47, 2).Value = "POSITION_OPENED"
$ws.Cells.Item(47, 3).Value = "SUI"
$ws.Cells.Item(47, 4).Value = "UNKNOWN"
$ws.Cells.Item(47, 5).Value = 3.348909889983477
$ws.Cells.Item(47, 6).Value = 2400
$ws.Cells.Item(47, 7).Value = 10
$ws.Cells.Item(47, 8).Value = 0.6539174731116587
$ws.Cells.Item(47, 11).Value = "SUCCESS"
